# Refresh crypto Price (D) and Volume(1h) (E) snapshot values to match the
# "Updated symbol list on Fri Jan  6 13:45:40 UTC 2023 with GitHub Actions" run.
# Cells are plain text in this sheet, so NumberFormat is forced to "@" (Text)
# before assigning, which keeps Excel/COM from reinterpreting the strings as
# numbers/percentages and rounding or reformatting them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '256.17'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.18%'

# Row 3: OKB
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.45'

# Row 4: HuobiToken
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.664'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.23%'

# Row 5: Cronos
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05925'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.23%'

# Row 6: KuCoinToken
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.604'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.56%'

# Row 7: MXToken
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8523'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-1.94%'

# Row 8: FTXToken
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9102'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.63%'

# Row 9: WazirX
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-2.04%'

# Row 10: LiechtensteinCryptoassetsExchange
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.04201'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '12.83%'

# Row 11: MandalaExchangeToken
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.57%'

# Row 12: BitrueCoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03028'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-5.42%'

# Row 13: BitMartToken
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09093'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.82%'

# Row 14: BitForexToken
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001524'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.39%'

# Row 15: One
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006026'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.33%'

# Row 16: TigerCash
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006024'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.61%'

# Row 17: LEO
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.469'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.24%'

# Row 18: GateToken
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.41%'

# Row 19: BTSEToken
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.90%'

# Row 21: ProBitToken
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1283'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.06%'

# Row 22: MCDex
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.854'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.18%'

# Row 23: CoinExToken
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04193'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.01%'

# Row 24: BitKan
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001215'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.39%'

# Row 25: HotbitToken
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004692'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '9.37%'

# Row 26: NitroEx
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.04%'

# Row 27: UpBots
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001522'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '1.30%'

# Row 40: IDEX
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03791'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.80%'

# Row 41: KickToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006257'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '58.40%'

# Row 42: BKEXToken
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.53%'

# Row 43: CEJI
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002311'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.51%'

# Row 44: LocalTraders
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01452'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '26.68%'

# Row 45: CoinLion
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005131'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-6.81%'

# Row 46: Kangarootoken
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.16%'

# Row 47: CoinbaseStockToken
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.04997'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-37.99%'

# Row 48: BOLO
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '10,458.58%'

# Row 49: CryptobidCoin
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.16%'

# Row 50: SpecialPowerGold
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.16%'
